$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update values in table 1 (Arquitetura Alpha / Camadas ocultas)
$ws.Range("I3").Value = 60
$ws.Range("I5").Value = 50

# Add new "Camada 4" row with merged cells, mirroring existing rows' layout
$ws.Range("H6").Value = "Camada 4"
$ws.Range("I6:J6").Merge()
$ws.Range("I6").Value = 80

# Fix typo in the model label that appears in both tables
$ws.Range("A2").Value = "Modelo de Validação Cruzada"
$ws.Range("A11").Value = "Modelo de Validação Cruzada"

$ws.Range("L19").Select()
